$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @('Cluster name', 'Active cases'),
  @('574 Plummer Street Building B1 & B2 Port Melbourne', 10),
  @('A1 Bakery Brunswick', 6),
  @('Acquire BPO Southbank', 36),
  @('Al Haj Halal Meats Glenroy', 71),
  @('Al-Taqwa College Truganina', 28),
  @('Amiga Montessori Craigieburn', 22),
  @('Broadmeadows Medical Centre Broadmeadows', 5),
  @('CMV Truck and Bus Dealership Springvale', 5),
  @('CS Square Caroline Springs', 5),
  @('Campbellfield Heights Primary School', 6),
  @('Can Panel Cambellfield', 6),
  @('Cedars Medical Clinic Coburg', 44),
  @('City of Greater Shepparton Community', 6),
  @('City of Hobsons Bay Community', 12),
  @('City of Wyndham Community', 7),
  @('Classy Cabinets and Kitchens Craigieburn', 16),
  @('Coles Barkly Square Brunswick August', 5),
  @('Coles Broadmeadows Central Shopping Centre', 11),
  @('Coles Campbellfield Plaza Campbellfield', 12),
  @('Coles Coburg North Village August', 8),
  @('Coles Greenvale Shopping Centre', 6),
  @('Coles Pakenham Place Shopping Centre', 6),
  @('Coles Roxburgh Village Roxburgh Park', 17),
  @('Community Kids Meadow Heights', 22),
  @('Concept Caravans Campbellfield', 5),
  @('Costco Wholesale Epping', 18),
  @('Crossroads Logistics Sunshine North', 6),
  @('Croydon Orthodontics', 8),
  @('Elite Smart Community Care Campbellfield', 5),
  @('Fitzroy Community School Fitzroy North', 50),
  @('Glenroy West Primary School', 7),
  @('Gloria Jeans Coffees Broadmeadows Central', 5),
  @('Health Care Providers Association South Melbourne', 12),
  @('IGA Meadow Heights Shopping Centre Meadow Heights', 7),
  @('Ilim College Glenroy Campus Hadfield', 25),
  @('Ilim Learning Sanctuary Glenroy', 14),
  @('Industrial Galvanizers Valmont Coatings Campbellfield', 13),
  @('Islamic College of Melbourne Tarneit', 7),
  @('KFC Fawkner', 21),
  @('Kasr Sweets Coolaroo', 22),
  @('Kids House Early Learning Cheltenham', 10),
  @('Learning Nest Early Learning Centre Meadow Heights', 10),
  @('Louis Vuitton Melbourne Crown Southbank', 5),
  @('Malvern Health and Fitness Clinic Malvern', 7),
  @('Melbourne Metropolitan Remand Centre Ravenhall', 5),
  @('Melbourne Truck Repairs Campbellfield', 5),
  @('Mercy Hospital for Women Heidelberg', 6),
  @('Montessori Beginnings Greenvale', 6),
  @('MyCentre Childcare Broadmeadows', 19),
  @('Newbury Child and Community Centre Craigieburn', 10),
  @('Newport Football Club Altona North', 7),
  @('Nido Early School Moonee Ponds', 18),
  @('Nino Early Learning Adventures Lalor', 5),
  @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 42),
  @('OnQ Plumbing and Excavations Craigieburn', 8),
  @('Oporto Coolaroo', 14),
  @('Paisley Park Early Learning Centre Bundoora', 8),
  @('Panorama Construction Site Whitehorse Rd Box Hill', 67),
  @('Pascoe Vale Primary School Pascoe Vale', 6),
  @('Pearl Street Child Care Centre Glenroy', 11),
  @('People First Healthcare Home Residence Disability Support Taylors Lakes', 5),
  @('Serco Mill Park', 9),
  @('Southern Cross Station Crew Room Tissue Box Docklands', 5),
  @('St Vincents Hospital Emergency Department Melbourne', 7),
  @('Tek Foods Somerton', 9),
  @('The Homestead Child and Family Centre Roxburgh Park', 17),
  @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 5),
  @('Tip Top Warehouse Dandenong', 10),
  @('Total Window Concepts Hoppers Crossing', 7),
  @('Tunstall Fresh Tunstall Square Shopping Centre Doncaster East', 5),
  @('Victorian Civil and Administrative Tribunal Melbourne', 5),
  @('Werribee Mercy Hospital Emergency Department', 7),
  @('Western Health Footscray Hospital Emergency Department', 6),
  @('Woolworths Coburg Station', 5),
  @('Woolworths Greenvale Lakes Roxburgh Park', 15),
  @('Yara Childcare Centre Truganina', 30)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
}

# Drop the now-obsolete 'Unnamed: 2' / 'Unnamed: 3' columns (delete right-to-left)
$ws.Columns.Item(4).Delete()
$ws.Columns.Item(3).Delete()
